# Auto-generated Excel COM-interop script
# Applies numeric/text corrections to the 广州-漫展信息 workbook
# as described in the commit "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# ---- 展览 (sheet1) ----
$ws1.Cells.Item(2, 6).Value = 299  # F2: 297 -> 299
$ws1.Cells.Item(3, 6).Value = 480  # F3: 478 -> 480
$ws1.Cells.Item(4, 6).Value = 218  # F4: 217 -> 218
$ws1.Cells.Item(4, 7).Value = "不可售"  # G4: 68.8 -> 不可售
$ws1.Cells.Item(5, 7).Value = "不可售"  # G5: 68 -> 不可售
$ws1.Cells.Item(6, 6).Value = 292  # F6: 291 -> 292
$ws1.Cells.Item(7, 6).Value = 7350  # F7: 7326 -> 7350
$ws1.Cells.Item(8, 6).Value = 85  # F8: 84 -> 85
$ws1.Cells.Item(9, 6).Value = 66  # F9: 65 -> 66
$ws1.Cells.Item(10, 6).Value = 3037  # F10: 2852 -> 3037
$ws1.Cells.Item(12, 6).Value = 576  # F12: 575 -> 576
$ws1.Cells.Item(13, 6).Value = 586  # F13: 583 -> 586
$ws1.Cells.Item(17, 6).Value = 744  # F17: 743 -> 744
$ws1.Cells.Item(19, 6).Value = 54  # F19: 53 -> 54
$ws1.Cells.Item(20, 6).Value = 195  # F20: 186 -> 195
$ws1.Cells.Item(22, 6).Value = 222  # F22: 218 -> 222
$ws1.Cells.Item(23, 6).Value = 120  # F23: 119 -> 120
$ws1.Cells.Item(24, 6).Value = 360  # F24: 358 -> 360
$ws1.Cells.Item(26, 6).Value = 1067  # F26: 1065 -> 1067
$ws1.Cells.Item(27, 6).Value = 71  # F27: 70 -> 71
$ws1.Cells.Item(28, 6).Value = 112  # F28: 106 -> 112
$ws1.Cells.Item(29, 6).Value = 2093  # F29: 2089 -> 2093
$ws1.Cells.Item(30, 6).Value = 613  # F30: 606 -> 613
$ws1.Cells.Item(31, 6).Value = 23  # F31: 22 -> 23
$ws1.Cells.Item(32, 6).Value = 24  # F32: 21 -> 24
$ws1.Cells.Item(34, 6).Value = 572  # F34: 571 -> 572
$ws1.Cells.Item(35, 6).Value = 28  # F35: 25 -> 28

# ---- 演出 (sheet2) ----
$ws2.Cells.Item(4, 6).Value = 301  # F4: 299 -> 301
$ws2.Cells.Item(5, 6).Value = 312  # F5: 311 -> 312

# ---- 本地生活 (sheet3) ----
$ws3.Cells.Item(2, 6).Value = 402  # F2: 399 -> 402

# ---- 全部类型 (sheet4) ----
$ws4.Cells.Item(2, 6).Value = 402  # F2: 399 -> 402
$ws4.Cells.Item(3, 6).Value = 299  # F3: 297 -> 299
$ws4.Cells.Item(4, 6).Value = 480  # F4: 478 -> 480
$ws4.Cells.Item(5, 6).Value = 218  # F5: 217 -> 218
$ws4.Cells.Item(5, 7).Value = "不可售"  # G5: 68.8 -> 不可售
$ws4.Cells.Item(6, 7).Value = "不可售"  # G6: 68 -> 不可售
$ws4.Cells.Item(7, 6).Value = 292  # F7: 291 -> 292
$ws4.Cells.Item(8, 6).Value = 7350  # F8: 7327 -> 7350
$ws4.Cells.Item(9, 6).Value = 85  # F9: 84 -> 85
$ws4.Cells.Item(10, 6).Value = 66  # F10: 65 -> 66
$ws4.Cells.Item(12, 6).Value = 3037  # F12: 2853 -> 3037
$ws4.Cells.Item(14, 6).Value = 576  # F14: 575 -> 576
$ws4.Cells.Item(15, 6).Value = 586  # F15: 583 -> 586
$ws4.Cells.Item(20, 6).Value = 301  # F20: 299 -> 301
$ws4.Cells.Item(21, 6).Value = 312  # F21: 311 -> 312
$ws4.Cells.Item(23, 6).Value = 744  # F23: 743 -> 744
$ws4.Cells.Item(25, 6).Value = 54  # F25: 53 -> 54
$ws4.Cells.Item(26, 6).Value = 195  # F26: 186 -> 195
$ws4.Cells.Item(31, 6).Value = 222  # F31: 218 -> 222
$ws4.Cells.Item(32, 6).Value = 120  # F32: 119 -> 120
$ws4.Cells.Item(33, 6).Value = 360  # F33: 358 -> 360
$ws4.Cells.Item(35, 6).Value = 1067  # F35: 1065 -> 1067
$ws4.Cells.Item(36, 6).Value = 71  # F36: 70 -> 71
$ws4.Cells.Item(37, 6).Value = 112  # F37: 106 -> 112
$ws4.Cells.Item(38, 6).Value = 2093  # F38: 2089 -> 2093
$ws4.Cells.Item(39, 6).Value = 613  # F39: 606 -> 613
$ws4.Cells.Item(40, 6).Value = 23  # F40: 22 -> 23
$ws4.Cells.Item(41, 6).Value = 24  # F41: 21 -> 24
$ws4.Cells.Item(43, 6).Value = 572  # F43: 571 -> 572
$ws4.Cells.Item(44, 6).Value = 28  # F44: 25 -> 28
